# refactor: the logic class structure in xml files
# Insert a new attribute-flag row ("Force") above the existing "Upload" row
# on the Property1 sheet, shifting the header/data rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 8 (pushes old row 8 "Upload" -> row 9, etc.)
$ws.Rows.Item(8).Insert()

# Copy formatting from the row below (the shifted-down "Upload" row) so the
# new row matches the existing boolean-flag row styling (label cell style +
# TRUE/FALSE cell style across B:S).
$ws.Range("A9:S9").Copy()
$ws.Range("A8:S8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values: label "Force" plus all-FALSE flags.
$ws.Range("A8").Value = "Force"
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
foreach ($col in $cols) {
    $ws.Range($col + "8").Value = $false
}

# Restore the frozen-pane split to include the newly inserted row, and move
# the active cell in the lower pane to A9 (mirrors the author's recorded
# cursor position after the edit).
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A11").Select()
$win.FreezePanes = $true
$ws.Range("A9").Select()
